$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 2

# B3:B10 share one formula (B1+B2 filled down) -> becomes a shared formula
# group (t="shared") spanning B3:B10, yielding the Fibonacci sequence.
$ws.Range("B3:B10").Formula = "=B1+B2"

[void]$ws.Range("B3:B10").Select()
